$d = $word.ActiveDocument

# --------------------------------------------------------------------------
# Locate the bullet paragraph ending "... you run it with the python3
# command." -- the last item of the "For the next team" bullet list before
# this edit. The two new bullets get appended right after it.
# --------------------------------------------------------------------------
$anchorIndex = -1
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $cand = $d.Paragraphs.Item($i)
    if ($cand.Range.Text -like "*python3 command.*") {
        $anchorIndex = $i
        break
    }
}

if ($anchorIndex -eq -1) {
    throw "Could not find anchor paragraph ending in 'python3 command.'"
}

# A minimal OOXML package wrapper so InsertXML can splice a fully-formed
# <w:p> (with its own <w:proofErr/> spell/grammar markers) into the body.
$pkgOpen = @'
<?xml version="1.0" standalone="yes"?>
<?mso-application progid="Word.Document"?>
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">
<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">
<pkg:xmlData>
<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
<w:body>
'@

$pkgClose = @'
</w:body>
</w:document>
</pkg:xmlData>
</pkg:part>
</pkg:package>
'@

# --- New bullet 1: VNC credentials ---------------------------------------
$anchor = $d.Paragraphs.Item($anchorIndex)
$r = $anchor.Range
$r.Collapse(0)
$r.InsertParagraphAfter()
$newPara1 = $d.Paragraphs.Item($anchorIndex + 1)

$body1 = @'
<w:p>
<w:pPr>
<w:pStyle w:val="ListParagraph"/>
<w:numPr>
<w:ilvl w:val="0"/>
<w:numId w:val="1"/>
</w:numPr>
</w:pPr>
<w:r>
<w:t xml:space="preserve">VNC username is pi, password is </w:t>
</w:r>
<w:proofErr w:type="spellStart"/>
<w:r>
<w:t>nightswatch</w:t>
</w:r>
<w:proofErr w:type="spellEnd"/>
</w:p>
'@

$newPara1.Range.InsertXML($pkgOpen + $body1 + $pkgClose)

# --- New bullet 2: repo / running-game copy location ----------------------
$r2 = $d.Paragraphs.Item($anchorIndex + 1).Range
$r2.Collapse(0)
$r2.InsertParagraphAfter()
$newPara2 = $d.Paragraphs.Item($anchorIndex + 2)

$body2 = @'
<w:p>
<w:pPr>
<w:pStyle w:val="ListParagraph"/>
<w:numPr>
<w:ilvl w:val="0"/>
<w:numId w:val="1"/>
</w:numPr>
</w:pPr>
<w:r>
<w:t xml:space="preserve">Copy of repo is on desktop. Copy of game that is </w:t>
</w:r>
<w:proofErr w:type="gramStart"/>
<w:r>
<w:t>actually running</w:t>
</w:r>
<w:proofErr w:type="gramEnd"/>
<w:r>
<w:t xml:space="preserve"> is in </w:t>
</w:r>
<w:r>
<w:t>/home/pi/Code</w:t>
</w:r>
</w:p>
'@

$newPara2.Range.InsertXML($pkgOpen + $body2 + $pkgClose)
